# Add an "active" column as the new first column on the "Debts" and
# "Fixed Assets" sheets (commit message: "add active column to debts and
# fixed assets").

$wb = $excel.ActiveWorkbook

# --- Debts sheet: insert a new column A, label it "active" -----------------
$wsDebts = $wb.Worksheets.Item("Debts")
$wsDebts.Columns("A:A").Insert()
$wsDebts.Range("A1").Value = "active"

# --- Fixed Assets sheet: insert a new column A, label it "active" ----------
$wsAssets = $wb.Worksheets.Item("Fixed Assets")
$wsAssets.Columns("A:A").Insert()
$wsAssets.Range("A1").Value = "active"

# --- Restore/update the on-screen selections ---------------------------
# Fixed Assets ends up showing E8 selected (no longer the active tab)...
$wsAssets.Activate()
$wsAssets.Range("E8").Select() | Out-Null

# ...while Debts becomes the active tab, with C14 selected.
$wsDebts.Activate()
$wsDebts.Range("C14").Select() | Out-Null
